$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()

$ws.Range("H4").Value = 302.45456
$ws.Range("I4").Value = 312.7
$ws.Range("J4").Value = 200
$ws.Range("K4").Value = 312.7
$ws.Range("L4").Value = 200
$ws.Range("M4").Value = -198.7
$ws.Range("N4").Value = -428

$ws.Range("H9").Value = 375
$ws.Range("I9").Value = 250
$ws.Range("J9").Value = 500
$ws.Range("K9").Value = 250
$ws.Range("L9").Value = 500
$ws.Range("M9").Value = -81
$ws.Range("N9").Value = -838

$ws.Range("H12").Value = 136.85715
$ws.Range("I12").Value = 149.5
$ws.Range("K12").Value = 149.5
$ws.Range("M12").Value = 20.5

$ws.Range("H33").Value = 183.6
$ws.Range("I33").Value = 98.55556
$ws.Range("K33").Value = 98.55556
$ws.Range("M33").Value = 130.44444

$ws.Range("H116").Value = 2300
$ws.Range("I116").Value = 2300
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2300
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1142
$ws.Range("N116").ClearContents()

$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("M125").ClearContents()
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 858.3
$ws.Range("I32").Value = 810.8889
$ws.Range("K32").Value = 810.8889
$ws.Range("M32").Value = -523.8889

$ws.Range("H45").Value = 1456
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()

$ws.Range("H102").Value = 2334.6667
$ws.Range("I102").Value = 2334.6667
$ws.Range("K102").Value = 2334.6667
$ws.Range("M102").Value = -712.6667000000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 24681.4
$ws.Range("I26").Value = 24681.4
$ws.Range("K26").Value = 24681.4
$ws.Range("M26").Value = -24389.4

$ws.Range("H99").Value = 3126.8
$ws.Range("I99").Value = 3250.1428
$ws.Range("J99").Value = 1400
$ws.Range("K99").Value = 3250.1428
$ws.Range("L99").Value = 1400
$ws.Range("M99").Value = -1752.1428
$ws.Range("N99").Value = -4396

$ws.Range("H105").Value = 1652
$ws.Range("I105").Value = 2012.25
$ws.Range("J105").Value = 211
$ws.Range("K105").Value = 2012.25
$ws.Range("L105").Value = 211
$ws.Range("M105").Value = -265.25
$ws.Range("N105").Value = -3705

$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()

$ws.Range("H99").Value = 5000000
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()

$ws.Range("H109").Value = 71000
$ws.Range("J109").Value = 71000
$ws.Range("L109").Value = 71000
$ws.Range("N109").Value = -73080

$ws.Range("H126").Value = 5000000
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws.Range("H134").Value = 1242.2222
$ws.Range("I134").Value = 1270.75
$ws.Range("J134").Value = 1014
$ws.Range("K134").Value = 3812.25
$ws.Range("L134").Value = 3042
$ws.Range("M134").Value = -1277.25
$ws.Range("N134").Value = -8112

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 5730
$ws.Range("I14").Value = 5730
$ws.Range("K14").Value = 17190
$ws.Range("M14").Value = -17017

$ws.Range("H17").Value = 448.57144
$ws.Range("I17").Value = 245
$ws.Range("J17").Value = 530
$ws.Range("K17").Value = 735
$ws.Range("L17").Value = 1590
$ws.Range("M17").Value = -566
$ws.Range("N17").Value = -1928

$ws.Range("H26").Value = 100
$ws.Range("I26").Value = 100
$ws.Range("J26").Value = 100
$ws.Range("K26").Value = 300
$ws.Range("L26").Value = 300
$ws.Range("M26").Value = -12
$ws.Range("N26").Value = -876

$ws.Range("H34").Value = 1601.4
$ws.Range("J34").Value = 4999.6665
$ws.Range("L34").Value = 14998.9995
$ws.Range("N34").Value = -15166.9995

$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws.Range("H100").Value = 8000
$ws.Range("J100").Value = 8000
$ws.Range("L100").Value = 24000
$ws.Range("N100").Value = -25622

$ws.Range("H114").Value = 1750
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()

$ws.Range("H121").Value = 891.5
$ws.Range("J121").Value = 891.5
$ws.Range("L121").Value = 2674.5
$ws.Range("N121").Value = -5294.5

$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("M129").ClearContents()

$ws.Range("H130").Value = 2130
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 2130
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 6390
$ws.Range("N130").Value = -16430
$ws.Range("M130").ClearContents()

$ws.Range("H131").Value = 5000
$ws.Range("J131").Value = 5000
$ws.Range("L131").Value = 15000
$ws.Range("N131").Value = -25080

$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("M134").ClearContents()

$ws.Range("H138").Value = 4009.6667
$ws.Range("I138").Value = 4009.6667
$ws.Range("K138").Value = 12029.0001
$ws.Range("M138").Value = -6889.000100000001

$ws.Range("H139").Value = 449.25
$ws.Range("I139").Value = 515.6667
$ws.Range("J139").Value = 250
$ws.Range("K139").Value = 1547.0001
$ws.Range("L139").Value = 750
$ws.Range("M139").Value = 3592.9999
$ws.Range("N139").Value = -11030

$ws.Range("H140").Value = 209.5
$ws.Range("I140").Value = 209.5
$ws.Range("K140").Value = 628.5
$ws.Range("M140").Value = 4551.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1700
$ws.Range("I7").Value = 1400
$ws.Range("K7").Value = 1400
$ws.Range("M7").Value = -1288

$ws.Range("H22").Value = 4499.1113
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 4499.1113
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 4499.1113
$ws.Range("N22").Value = -5089.1113
$ws.Range("M22").ClearContents()

$ws.Range("H27").Value = 4499.1113
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 4499.1113
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 4499.1113
$ws.Range("N27").Value = -4713.1113
$ws.Range("M27").ClearContents()

$ws.Range("H40").Value = 650675.6
$ws.Range("I40").Value = 22580
$ws.Range("J40").Value = 1697501.6
$ws.Range("K40").Value = 22580
$ws.Range("L40").Value = 1697501.6
$ws.Range("M40").Value = -22444
$ws.Range("N40").Value = -1697773.6

$ws.Range("H46").Value = 254975.12
$ws.Range("I46").Value = 1001200.5
$ws.Range("J46").Value = 6233.3335
$ws.Range("K46").Value = 1001200.5
$ws.Range("L46").Value = 6233.3335
$ws.Range("M46").Value = -1001012.5
$ws.Range("N46").Value = -6609.3335

$ws.Range("H98").Value = 57177.5
$ws.Range("J98").Value = 57177.5
$ws.Range("L98").Value = 57177.5
$ws.Range("N98").Value = -63167.5

$ws.Range("H100").Value = 2999.5
$ws.Range("I100").Value = 2999.5
$ws.Range("K100").Value = 2999.5
$ws.Range("M100").Value = -2458.5

$ws.Range("H126").Value = 1700
$ws.Range("I126").Value = 1400
$ws.Range("K126").Value = 4200
$ws.Range("M126").Value = -1730

$ws.Range("H132").Value = 27250
$ws.Range("I132").Value = 27250
$ws.Range("K132").Value = 81750
$ws.Range("M132").Value = -79220

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 52625
$ws.Range("I29").Value = 43500
$ws.Range("J29").Value = 80000
$ws.Range("K29").Value = 43500
$ws.Range("L29").Value = 80000
$ws.Range("M29").Value = -43210
$ws.Range("N29").Value = -80580

$ws.Range("H125").Value = 39215
$ws.Range("J125").Value = 39215
$ws.Range("L125").Value = 39215
$ws.Range("N125").Value = -49055

$ws.Range("H126").Value = 5115.4
$ws.Range("I126").Value = 3526
$ws.Range("K126").Value = 10578
$ws.Range("M126").Value = -8108

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()
